# Generate Report for Handoff
# Adds a new handed-off file (e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md) as a
# new row across the Overview / zh-cn / de-de sheets + their tables.

$wb = $excel.ActiveWorkbook

$commitSha = "a6b295adb4f1680d43cd53b3aad9338a7f879817"
$newFile   = "e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md"
$hlTarget  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newFile"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Set-PlainText($range, $text) {
    # Force literal text (not auto-converted to bool/number/date) by using a
    # leading apostrophe, then reset the style so the quote-prefix marker
    # doesn't leak into the saved style.
    if ($text -eq "") {
        $range.Value2 = "'"
    } else {
        $range.Value2 = "'" + $text
    }
    $range.Style = "Normal"
}

function Add-FileHyperlink($ws, $range, $target, $display) {
    $ws.Hyperlinks.Add(
        $range,
        $target,
        [System.Type]::Missing,
        [System.Type]::Missing,
        $display
    ) | Out-Null
    # Match the workbook's existing custom "HyperLink" cell style (blue
    # FF6495ED, single underline) instead of Excel's themed default.
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# ============================= Overview sheet =============================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

Set-PlainText $wsOverview.Range("A3") $newFile
Set-PlainText $wsOverview.Range("C3") ".md"
Set-PlainText $wsOverview.Range("D3") ""
Set-PlainText $wsOverview.Range("E3") "Ready for handoff"
Set-PlainText $wsOverview.Range("F3") "Ready for handoff"

$wsOverview.Range("G3").Value2 = "2016-08-20 12:43:46"
$wsOverview.Range("G3").NumberFormat = $dateFmt

Add-FileHyperlink $wsOverview $wsOverview.Range("B3") $hlTarget "e2e\$newFile"

# =============================== zh-cn sheet ================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

Set-PlainText $wsZhCn.Range("B3") ".md"
Set-PlainText $wsZhCn.Range("C3") "Ready for handoff"
Set-PlainText $wsZhCn.Range("D3") "e2e"
Set-PlainText $wsZhCn.Range("E3") "ht"
Set-PlainText $wsZhCn.Range("F3") "False"
Set-PlainText $wsZhCn.Range("G3") "e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.df51c59b9feb0886828735e027751f7265be5dc1.zh-cn.xlf"
Set-PlainText $wsZhCn.Range("I3") ""
Set-PlainText $wsZhCn.Range("J3") ""
Set-PlainText $wsZhCn.Range("L3") ""
Set-PlainText $wsZhCn.Range("M3") "True"
Set-PlainText $wsZhCn.Range("N3") ""
Set-PlainText $wsZhCn.Range("O3") "False"
Set-PlainText $wsZhCn.Range("P3") ""

$wsZhCn.Range("H3").Value2 = "2016-08-20 12:43:42"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFmt

Add-FileHyperlink $wsZhCn $wsZhCn.Range("A3") $hlTarget $newFile

# =============================== de-de sheet ================================
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

Set-PlainText $wsDeDe.Range("B3") ".md"
Set-PlainText $wsDeDe.Range("C3") "Ready for handoff"
Set-PlainText $wsDeDe.Range("D3") "e2e"
Set-PlainText $wsDeDe.Range("E3") "ht"
Set-PlainText $wsDeDe.Range("F3") "False"
Set-PlainText $wsDeDe.Range("G3") "e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.df51c59b9feb0886828735e027751f7265be5dc1.de-de.xlf"
Set-PlainText $wsDeDe.Range("I3") ""
Set-PlainText $wsDeDe.Range("J3") ""
Set-PlainText $wsDeDe.Range("L3") ""
Set-PlainText $wsDeDe.Range("M3") "True"
Set-PlainText $wsDeDe.Range("N3") ""
Set-PlainText $wsDeDe.Range("O3") "False"
Set-PlainText $wsDeDe.Range("P3") ""

$wsDeDe.Range("H3").Value2 = "2016-08-20 12:43:46"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFmt

Add-FileHyperlink $wsDeDe $wsDeDe.Range("A3") $hlTarget $newFile

Write-Host "Handoff report row added to Overview, zh-cn, de-de sheets."
